$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 and 11 swap: OKB <-> Avalanche (content moves, rank index A stays)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "30.80"
$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "51.34"
$ws.Range("E11").Value = "  +7.65%  "

# Price / Volume(1h) updates for remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.041.52"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.214.13"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.74"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.30"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0780"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.37"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.559.58"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.81"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.177.10"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.733"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.977.33"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.19"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.75"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.50"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.09"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.81"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.04"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.28"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.51"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.55"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.05"
$ws.Range("E34").Value = "  +6.16%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0711"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.55"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.071.29"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.28"
$ws.Range("E44").Value = "  +10.29%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  -10.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.434.72"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.45"
$ws.Range("E51").Value = "  +0.22%  "
